# Apply updated model predictions for Sheet1 ("Train Results") and
# Sheet2 ("Test Results"): existing rows get new feature/prediction values,
# and additional rows are appended for both train and test data.

$wb = $excel.ActiveWorkbook

# ---- Sheet1: "Train Results" ----
$ws1 = $wb.Worksheets.Item("Train Results")

$rows1 = New-Object 'object[]' 43
$rows1[0] = @(0,40,4,0,28,4,4,20,2.31,2.346740961074829)
$rows1[1] = @(4,0,8,24,4,16,44,0,2.98,3.042011737823486)
$rows1[2] = @(4,16,0,20,4,12,44,0,3.22,3.093508005142212)
$rows1[3] = @(0,24,4,0,8,4,36,24,2.45,2.33285665512085)
$rows1[4] = @(4,0,8,4,24,12,44,4,3.17,3.00151801109314)
$rows1[5] = @(0,16,8,4,16,20,36,0,3.07,2.908205032348633)
$rows1[6] = @(4,20,4,4,16,0,52,0,3.35,3.081407070159912)
$rows1[7] = @(4,0,12,4,4,16,52,8,2.81,2.775321483612061)
$rows1[8] = @(4,12,8,0,8,8,56.00000000000001,4,2.88,2.888345956802368)
$rows1[9] = @(4,4,4,8,28,16,32,4,3.04,3.069018363952637)
$rows1[10] = @(4,0,8,4,24,12,44,4,3.09,3.00151801109314)
$rows1[11] = @(4,8,0,12,20,8,48,0,3.11,3.116998434066772)
$rows1[12] = @(0,20,8,4,0,20,44,4,2.82,2.726865291595459)
$rows1[13] = @(0,20,8,4,0,20,44,4,2.74,2.726865291595459)
$rows1[14] = @(4,12,0,0,16,8,52,8,2.62,2.8892502784729)
$rows1[15] = @(4,12,8,4,16,8,48,0,3.19,3.022479772567749)
$rows1[16] = @(0,20,4,0,4,4,48,20,2.52,2.349213600158691)
$rows1[17] = @(4,8,4,4,16,12,48,4,2.94,2.976077795028687)
$rows1[18] = @(4,12,0,12,20,8,44,0,3.29,3.124362945556641)
$rows1[19] = @(4,12,0,12,20,8,44,0,3.25,3.124362945556641)
$rows1[20] = @(0,12,8,4,16,20,36,4,3.01,2.838096380233765)
$rows1[21] = @(4,12,4,4,20,16,32,8,2.96,2.951002836227417)
$rows1[22] = @(4,8,0,0,4,8,52,24,2.44,2.41313910484314)
$rows1[23] = @(4,0,4,0,20,8,52,12,2.92,2.730660200119019)
$rows1[24] = @(4,8,12,4,4,24,39.99999999999999,4,2.9,2.879618406295776)
$rows1[25] = @(4,0,12,16,4,12,52,0,2.96,3.003385543823242)
$rows1[26] = @(0,16,8,4,16,20,36,0,3.02,2.908205032348633)
$rows1[27] = @(4,12,8,0,8,8,56.00000000000001,4,2.79,2.888345956802368)
$rows1[28] = @(4,0,4,4,16,12,56.00000000000001,4,2.85,2.96811056137085)
$rows1[29] = @(0,4,8,4,16,20,44,4,2.94,2.8643798828125)
$rows1[30] = @(0,12,4,0,4,16,39.99999999999999,24,2.51,2.389402866363525)
$rows1[31] = @(0,0,8,4,16,12,52,8,2.99,2.737751007080078)
$rows1[32] = @(0,16,0,4,20,20,28,12,3.38,2.746609449386597)
$rows1[33] = @(4,16,8,0,12,0,48,12,2.56,2.541861057281494)
$rows1[34] = @(0,0,8,4,16,12,52,8,2.82,2.737751007080078)
$rows1[35] = @(4,12,4,4,12,16,36,12,2.86,2.839879989624023)
$rows1[36] = @(0,24,8,0,12,16,32,8,2.93,2.635068893432617)
$rows1[37] = @(4,32,8,4,4,20,24,4,2.84,2.897067785263062)
$rows1[38] = @(0,24,8,0,12,16,32,8,2.94,2.635068893432617)
$rows1[39] = @(4,40,0,0,12,4,36,4,3.16,3.104430198669434)
$rows1[40] = @(0,16,0,4,20,20,28,12,2.72,2.746609449386597)
$rows1[41] = @(4,20,4,4,16,0,52,0,3.21,3.081407070159912)
$rows1[42] = @(4,12,0,0,16,8,52,8,2.56,2.8892502784729)

$data1 = New-Object 'object[,]' $rows1.Count,10
for ($r = 0; $r -lt $rows1.Count; $r++) {
  for ($c = 0; $c -lt 10; $c++) {
    $data1[$r,$c] = $rows1[$r][$c]
  }
}
$ws1.Range("A2:J44").Value = $data1

# ---- Sheet2: "Test Results" ----
$ws2 = $wb.Worksheets.Item("Test Results")

$rows2 = New-Object 'object[]' 9
$rows2[0] = @(0,8,0,4,4,28,36,20,2.56,2.639162540435791)
$rows2[1] = @(4,16,8,0,12,0,48,12,2.67,2.541861057281494)
$rows2[2] = @(4,16,0,4,12,8,52,4,2.96,2.986655235290527)
$rows2[3] = @(4,8,4,4,16,12,48,4,2.92,2.976077795028687)
$rows2[4] = @(4,0,4,0,4,4,60,24,2.45,2.386898040771484)
$rows2[5] = @(4,0,4,0,20,8,52,12,2.64,2.730660200119019)
$rows2[6] = @(4,12,4,0,28,0,39.99999999999999,12,2.94,2.694723844528198)
$rows2[7] = @(4,4,12,0,0,4,52,24,2.22,2.366406440734863)
$rows2[8] = @(4,12,8,4,16,8,48,0,3.15,3.022479772567749)

$data2 = New-Object 'object[,]' $rows2.Count,10
for ($r = 0; $r -lt $rows2.Count; $r++) {
  for ($c = 0; $c -lt 10; $c++) {
    $data2[$r,$c] = $rows2[$r][$c]
  }
}
$ws2.Range("A2:J10").Value = $data2
